$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A91").Value = "2023-12-08 08:28:02"
$ws.Range("B91").Value = 0.0004

$ws.Range("A92").Value = "2023-12-08 08:28:24"
$ws.Range("B92").Value = 0.0006000000000000001

$ws.Range("A93").Value = "2023-12-08 08:29:31"
$ws.Range("B93").Value = 0.004600000000000001
